$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "ronenchen27@gmail.com"
$ws.Range("D6").Value = "danfogel100@gmail.com"
$ws.Range("F6").Value = "I love playing this game so much. Great entertainment and very funny to see the car falling to the river"

$ws.Range("F6").Select() | Out-Null

$wb.Save()
